# For every row in the "Recorded By" column (column G), if the value
# starts with "System, " move the leading "System" token to the end of
# the comma-separated list instead of the front.
#
# e.g. "System, dnasr281@gmail.com"          -> "dnasr281@gmail.com, System"
#      "System, system, backup@backdoor.com" -> "system, backup@backdoor.com, System"
#
# Values that do not start with "System, " (e.g. "dnasr281@gmail.com",
# "dnasr281@gmail.com, admin@admin.com") are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$recordedByCol = 7   # column G ("Recorded By")
$prefix = "System, "

$used = $ws.UsedRange
$lastRow = $used.Rows.Count

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, $recordedByCol)
    $text = $cell.Text

    if ($text -ne $null -and $text.StartsWith($prefix)) {
        $remainder = $text.Substring($prefix.Length)
        $cell.Value = $remainder + ", System"
    }
}
